$d = $word.ActiveDocument
$d.Paragraphs.LineSpacingRule = 0
